$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("caseType1")
$ws2 = $wb.Worksheets.Item("caseType1-vl")
$ws3 = $wb.Worksheets.Item("caseType2")

# "caseType1": bulk-upload template gains a new "Required" column, inserted
# right after "Case Property" (i.e. before the existing "Group" column).
$ws1.Columns("B").Insert() | Out-Null
$ws1.Range("B1").Value = "Required"
$ws1.Activate() | Out-Null
$ws1.Range("B:B").Select() | Out-Null

# "caseType1-vl" is untouched structurally; just leave the selection on
# column B (mirrors how the sibling sheets now highlight their B column).
$ws2.Activate() | Out-Null
$ws2.Range("B:B").Select() | Out-Null

# "caseType2" gets the same new "Required" column as "caseType1".
$ws3.Columns("B").Insert() | Out-Null
$ws3.Range("B1").Value = "Required"
$ws3.Activate() | Out-Null
$ws3.Range("B:B").Select() | Out-Null
